# Fruta / hortaliza, semanal
# Insert a new weekly record for "Poroto verde" (Terminal Hortofrutícola Agro
# Chillán) as row 57, pushing the existing rows 57..91 down to 58..92.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 57 - Excel shifts rows 57..91
# down to 58..92 and extends the used range/dimension accordingly, copying
# the formatting (including the date style on column D) from the row above.
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new data record.
$ws.Range("A57").Value = 7
$ws.Range("B57").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C57").Value = "Ñuble"
$ws.Range("D57").Value = 44790
$ws.Range("E57").Value = 16
$ws.Range("F57").Value = 100112031
$ws.Range("G57").Value = "Poroto verde"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 50
$ws.Range("K57").Value = 35000
$ws.Range("L57").Value = 35000
$ws.Range("M57").Value = 35000
$ws.Range("N57").Value = "`$/malla 25 kilos"
$ws.Range("O57").Value = "Región de Arica y Parinacota"
$ws.Range("P57").Value = 1400
$ws.Range("Q57").Value = 25
$ws.Range("R57").Value = "Hortaliza"
